$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cell, $value) {
    $cell.NumberFormat = "@"
    $cell.Value = $value
    $cell.Style = "Normal"
}

Set-TextValue $ws.Range("D2") '42.759.14'
Set-TextValue $ws.Range("E2") '  -0.10%  '
Set-TextValue $ws.Range("D3") '2.528.55'
Set-TextValue $ws.Range("E3") '  -1.44%  '
Set-TextValue $ws.Range("D4") '0.999'
Set-TextValue $ws.Range("E4") '  -0.10%  '
Set-TextValue $ws.Range("D5") '309.86'
Set-TextValue $ws.Range("E5") '  -0.55%  '
Set-TextValue $ws.Range("D6") '100.28'
Set-TextValue $ws.Range("E6") '  +2.08%  '
Set-TextValue $ws.Range("D7") '0.568'
Set-TextValue $ws.Range("E7") '  -0.90%  '
Set-TextValue $ws.Range("E8") '  +0.03%  '
Set-TextValue $ws.Range("E9") '  -1.84%  '
Set-TextValue $ws.Range("D10") '35.50'
Set-TextValue $ws.Range("E10") '  -0.27%  '
Set-TextValue $ws.Range("D11") '0.0806'
Set-TextValue $ws.Range("E11") '  -0.27%  '
Set-TextValue $ws.Range("D12") '7.32'
Set-TextValue $ws.Range("E12") '  -1.38%  '
Set-TextValue $ws.Range("E13") '  +1.24%  '
Set-TextValue $ws.Range("D14") '2.917.28'
Set-TextValue $ws.Range("E14") '  -1.68%  '
Set-TextValue $ws.Range("B15") 'Chainlink'
Set-TextValue $ws.Range("C15") 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link'
Set-TextValue $ws.Range("D15") '15.34'
Set-TextValue $ws.Range("E15") '  -3.25%  '
Set-TextValue $ws.Range("B16") 'WrappedEther'
Set-TextValue $ws.Range("C16") 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range("D16") '2.520.11'
Set-TextValue $ws.Range("E16") '  -1.30%  '
Set-TextValue $ws.Range("E17") '  -3.12%  '
Set-TextValue $ws.Range("D18") '42.732.62'
Set-TextValue $ws.Range("E18") '  -0.25%  '
Set-TextValue $ws.Range("E19") '  -0.84%  '
Set-TextValue $ws.Range("E20") '  -0.66%  '
Set-TextValue $ws.Range("D21") '12.28'
Set-TextValue $ws.Range("E21") '  -0.96%  '
Set-TextValue $ws.Range("D22") '69.33'
Set-TextValue $ws.Range("E22") '  -0.63%  '
Set-TextValue $ws.Range("D23") '243.15'
Set-TextValue $ws.Range("E23") '  -2.10%  '
Set-TextValue $ws.Range("E24") '  -2.20%  '
Set-TextValue $ws.Range("E25") '  -1.39%  '
Set-TextValue $ws.Range("E26") '  +0.10%  '
Set-TextValue $ws.Range("D27") '25.37'
Set-TextValue $ws.Range("E27") '  -6.01%  '
Set-TextValue $ws.Range("E28") '  -2.26%  '
Set-TextValue $ws.Range("D29") '10.17'
Set-TextValue $ws.Range("E29") '  -0.45%  '
Set-TextValue $ws.Range("D30") '38.53'
Set-TextValue $ws.Range("E30") '  -2.77%  '
Set-TextValue $ws.Range("D31") '160.49'
Set-TextValue $ws.Range("E31") '  +0.83%  '
Set-TextValue $ws.Range("E32") '  -0.01%  '
Set-TextValue $ws.Range("E33") '  +8.75%  '
Set-TextValue $ws.Range("E34") '  -0.06%  '
Set-TextValue $ws.Range("D35") '0.0785'
Set-TextValue $ws.Range("E35") '  -1.34%  '
Set-TextValue $ws.Range("D36") '18.44'
Set-TextValue $ws.Range("E36") '  -0.45%  '
Set-TextValue $ws.Range("E37") '  -5.58%  '
Set-TextValue $ws.Range("D38") '1.97'
Set-TextValue $ws.Range("E38") '  -6.50%  '
Set-TextValue $ws.Range("E39") '  -0.74%  '
Set-TextValue $ws.Range("E40") '  -0.25%  '
Set-TextValue $ws.Range("D41") '22.38'
Set-TextValue $ws.Range("E41") '  -1.80%  '
Set-TextValue $ws.Range("D42") '4.18'
Set-TextValue $ws.Range("E42") '  +1.70%  '
Set-TextValue $ws.Range("E43") '  +0.11%  '
Set-TextValue $ws.Range("E44") '  -0.18%  '
Set-TextValue $ws.Range("D45") '3.29'
Set-TextValue $ws.Range("E45") '  +2.91%  '
Set-TextValue $ws.Range("D46") '2.001.74'
Set-TextValue $ws.Range("E46") '  +0.24%  '
Set-TextValue $ws.Range("D47") '8.96'
Set-TextValue $ws.Range("E47") '  -0.38%  '
Set-TextValue $ws.Range("D48") '2.770.50'
Set-TextValue $ws.Range("E48") '  -1.65%  '
Set-TextValue $ws.Range("E49") '  -2.72%  '
Set-TextValue $ws.Range("D50") '79.27'
Set-TextValue $ws.Range("E50") '  -2.78%  '
Set-TextValue $ws.Range("D51") '71.88'
Set-TextValue $ws.Range("E51") '  -2.94%  '
